$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.899.37'
$ws.Range('E2').Value = '  -0.11%  '

# Row 3
$ws.Range('D3').Value = '3.409.62'
$ws.Range('E3').Value = '  -0.81%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').Value = '''409.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.04%  '

# Row 6
$ws.Range('D6').Value = '''129.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.23%  '

# Row 7
$ws.Range('D7').Value = '''0.621'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.25%  '

# Row 8
$ws.Range('D8').Value = '''1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.02%  '

# Row 9
$ws.Range('D9').Value = '''0.719'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.36%  '

# Row 10
$ws.Range('D10').Value = '''0.137'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.00%  '

# Row 11
$ws.Range('D11').Value = '''42.71'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.28%  '

# Row 12
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '3.951.31'
$ws.Range('E12').Value = '  -0.32%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = '''0.141'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.00%  '

# Row 14
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '''9.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.87%  '

# Row 15
$ws.Range('D15').Value = '''0.0000207'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.21%  '

# Row 16
$ws.Range('D16').Value = '''20.76'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.65%  '

# Row 17
$ws.Range('D17').Value = '3.391.76'
$ws.Range('E17').Value = '  +0.64%  '

# Row 18
$ws.Range('D18').Value = '''12.44'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.54%  '

# Row 19
$ws.Range('D19').Value = '''1.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.40%  '

# Row 20
$ws.Range('D20').Value = '61.816.14'
$ws.Range('E20').Value = '  -0.17%  '

# Row 21
$ws.Range('D21').Value = '''477.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +18.48%  '

# Row 22
$ws.Range('D22').Value = '''90.64'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.76%  '

# Row 23
$ws.Range('D23').Value = '''3.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.78%  '

# Row 24
$ws.Range('D24').Value = '''13.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.47%  '

# Row 25
$ws.Range('D25').Value = '''3.31'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.78%  '

# Row 26
$ws.Range('D26').Value = '''33.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.37%  '

# Row 27
$ws.Range('D27').Value = '''8.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.59%  '

# Row 28
$ws.Range('E28').Value = '  +0.35%  '

# Row 29
$ws.Range('D29').Value = '''7.68'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.01%  '

# Row 30
$ws.Range('D30').Value = '''2.75'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.03%  '

# Row 31
$ws.Range('D31').Value = '''11.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.15%  '

# Row 32
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = '''0.112'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.98%  '

# Row 33
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D33').Value = '''0.164'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.48%  '

# Row 34
$ws.Range('D34').Value = '''41.33'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.57%  '

# Row 35
$ws.Range('D35').Value = '''1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.04%  '

# Row 36
$ws.Range('D36').Value = '''56.49'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.24%  '

# Row 37
$ws.Range('D37').Value = '''0.0486'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.65%  '

# Row 38
$ws.Range('E38').Value = '  +0.09%  '

# Row 39
$ws.Range('D39').Value = '''149.63'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.25%  '

# Row 40
$ws.Range('D40').Value = '''3.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.79%  '

# Row 41
$ws.Range('E41').Value = '  +2.05%  '

# Row 42
$ws.Range('D42').Value = '''0.319'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.52%  '

# Row 43
$ws.Range('D43').Value = '''2.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.09%  '

# Row 44
$ws.Range('D44').Value = '''2.09'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.30%  '

# Row 45
$ws.Range('D45').Value = '''2.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.81%  '

# Row 46
$ws.Range('D46').Value = '''4.16'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.44%  '

# Row 47
$ws.Range('D47').Value = '''16.39'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.76%  '

# Row 48
$ws.Range('D48').Value = '''2.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +18.48%  '

# Row 49
$ws.Range('D49').Value = '''0.148'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +14.10%  '

# Row 50
$ws.Range('D50').Value = '''22.18'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.26%  '

# Row 51
$ws.Range('D51').Value = '''115.11'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +18.12%  '
